{"js": "// CIV-17609 updated GA documents to display main claim number\n//\n// Functional change: the label in front of the \"<<caseNumber>>\" merge\n// field reads \"Claim number:\" and should read \"Case number:\" instead\n// (the merge field itself, caseNumber, is unchanged). There is a second,\n// unrelated paragraph elsewhere in the template that already legitimately\n// reads \"Case number: <<claimNumber>>\" - that one must be left alone, so\n// we match on the old \"Claim number:\" text specifically (case-sensitive),\n// which is unique in the document.\n\nconst body = context.document.body;\n\nconst matches = body.search(\"Claim number:\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nmatches.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < matches.items.length; i++) {\n  matches.items[i].insertText(\"Case number:\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// The same commit also marks the built-in \"Default Paragraph Font\"\n// character style as semi-hidden (w:semiHidden) in styles.xml, matching\n// how \"Normal Table\" / \"No List\" / etc. are already flagged in this\n// template. Word's object model exposes this as Style.hidden.\ntry {\n  const styles = context.document.getStyles();\n  const defaultParaFont = styles.getByNameOrNullObject(\"Default Paragraph Font\");\n  await context.sync();\n\n  if (!defaultParaFont.isNullObject) {\n    defaultParaFont.hidden = true;\n    await context.sync();\n  }\n} catch (e) {\n  // Some hosts don't implement this legacy style flag - non-fatal.\n}\n", "ps1": "# CIV-17609 updated GA documents to display main claim number\n#\n# Functional change: the label in front of the \"<<caseNumber>>\" merge\n# field reads \"Claim number:\" and should read \"Case number:\" instead\n# (the merge field itself, caseNumber, is unchanged). There is a second,\n# unrelated paragraph elsewhere in the template that already legitimately\n# reads \"Case number: <<claimNumber>>\" - that one must be left alone, so\n# we match on the old \"Claim number:\" text specifically (case-sensitive),\n# which is unique in the document.\n\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"Claim number:\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Wrap = 1  # wdFindContinue\n$find.Forward = $true\n$find.Execute(\"Claim number:\", $true, $false, $false, $false, $false, $true, 1, $false, \"Case number:\", 2) | Out-Null\n\n# The same commit also marks the built-in \"Default Paragraph Font\"\n# character style as semi-hidden (w:semiHidden) in styles.xml, matching\n# how \"Normal Table\" / \"No List\" / etc. are already flagged in this\n# template. Word's object model exposes this as Style.Hidden.\ntry {\n    $style = $d.Styles(\"Default Paragraph Font\")\n    $style.Hidden = $true\n} catch {\n    # Some hosts don't implement this legacy style flag - non-fatal.\n}\n"}
